# Insert a new data row right before the current row 204, shifting the
# existing rows 204-276 down to 205-277 (formatting/styles move with them).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(204).Insert()

$row = 204

$ws.Cells.Item($row, 1).Value = 1                                       # A Mercado ID
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"      # B Mercado
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"                    # C Región
$ws.Cells.Item($row, 4).Value = 44784                                   # D Fecha
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 15                                      # E Codreg
$ws.Cells.Item($row, 6).Value = "Fruta"                                 # F Tipo
$ws.Cells.Item($row, 7).Value = 100108                                  # G Producto ID
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"            # H Producto
$ws.Cells.Item($row, 9).Value = 100108006                               # I Categoría ID
$ws.Cells.Item($row, 10).Value = "Plátano"                              # J Categoría
$ws.Cells.Item($row, 11).Value = "Sin especificar"                      # K Variedad
$ws.Cells.Item($row, 12).Value = "Pintón"                               # L Calidad
$ws.Cells.Item($row, 13).Value = 120                                    # M Volumen
$ws.Cells.Item($row, 14).Value = 21000                                  # N Precio mínimo
$ws.Cells.Item($row, 15).Value = 22000                                  # O Precio máximo
$ws.Cells.Item($row, 16).Value = 21500                                  # P Precio promedio ponderado
$ws.Cells.Item($row, 17).Value = "$/caja 20 kilos"                      # Q Unidad de comercialización
$ws.Cells.Item($row, 18).Value = "Ecuador"                              # R Origen
$ws.Cells.Item($row, 19).Value = 1075                                   # S Precio $/Kg
$ws.Cells.Item($row, 20).Value = 20                                     # T Kg / unidad
